# Update Work Week and Social Spending
# (Liberia GDP-per-Capita indicator workbook: refresh the "Data" values for the
# already-present years 1950-2008, and append the newly published years 2009-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Ensure the Data column keeps storing these as text (matches original t="s" cells)
$ws.Range("E2:E68").NumberFormat = "@"

# --- Update the existing 59 data rows (years 1950-2008) with the new GDP-per-capita-style values ---
$ws.Range("E2").Value = "3151"
$ws.Range("E3").Value = "3257"
$ws.Range("E4").Value = "3277"
$ws.Range("E5").Value = "3325"
$ws.Range("E6").Value = "3427"
$ws.Range("E7").Value = "3454"
$ws.Range("E8").Value = "3494"
$ws.Range("E9").Value = "3537"
$ws.Range("E10").Value = "3545"
$ws.Range("E11").Value = "3655"
$ws.Range("E12").Value = "3674"
$ws.Range("E13").Value = "3663"
$ws.Range("E14").Value = "3612"
$ws.Range("E15").Value = "3598"
$ws.Range("E16").Value = "3679"
$ws.Range("E17").Value = "3639"
$ws.Range("E18").Value = "4208"
$ws.Range("E19").Value = "4065"
$ws.Range("E20").Value = "4136"
$ws.Range("E21").Value = "4309"
$ws.Range("E22").Value = "4457"
$ws.Range("E23").Value = "4544"
$ws.Range("E24").Value = "4568"
$ws.Range("E25").Value = "4603"
$ws.Range("E26").Value = "4616"
$ws.Range("E27").Value = "4326"
$ws.Range("E28").Value = "4415"
$ws.Range("E29").Value = "4339"
$ws.Range("E30").Value = "4463"
$ws.Range("E31").Value = "4546"
$ws.Range("E32").Value = "4138"
$ws.Range("E33").Value = "3862"
$ws.Range("E34").Value = "3813"
$ws.Range("E35").Value = "3572"
$ws.Range("E36").Value = "3443"
$ws.Range("E37").Value = "3269"
$ws.Range("E38").Value = "3129"
$ws.Range("E39").Value = "3084"
$ws.Range("E40").Value = "3070"
$ws.Range("E41").Value = "3054"
$ws.Range("E42").Value = "1690"
$ws.Range("E43").Value = "1568.40432869831"
$ws.Range("E44").Value = "915.120585248132"
$ws.Range("E45").Value = "560.333014747494"
$ws.Range("E46").Value = "418.121582144975"
$ws.Range("E47").Value = "377.580086005337"
$ws.Range("E48").Value = "391.085917171772"
$ws.Range("E49").Value = "701.785385010882"
$ws.Range("E50").Value = "772.676468284599"
$ws.Range("E51").Value = "848.421168136595"
$ws.Range("E52").Value = "971.933016146115"
$ws.Range("E53").Value = "1097.24958690587"
$ws.Range("E54").Value = "1342.75876125062"
$ws.Range("E55").Value = "851.910950210634"
$ws.Range("E56").Value = "760.875776109511"
$ws.Range("E57").Value = "755.223591044875"
$ws.Range("E58").Value = "769.292953529305"
$ws.Range("E59").Value = "768.207040845143"
$ws.Range("E60").Value = "811.158033784576"

# --- Append the 8 new rows for years 2009-2016 ---
$ws.Range("A61").Value = 430
$ws.Range("B61").Value = "Liberia"
$ws.Range("C61").Value = "GDP per Capita"
$ws.Range("D61").Value = 2009
$ws.Range("E61").Value = "838.45203626212"
$ws.Range("A62").Value = 430
$ws.Range("B62").Value = "Liberia"
$ws.Range("C62").Value = "GDP per Capita"
$ws.Range("D62").Value = 2010
$ws.Range("E62").Value = "854.234972865737"
$ws.Range("A63").Value = 430
$ws.Range("B63").Value = "Liberia"
$ws.Range("C63").Value = "GDP per Capita"
$ws.Range("D63").Value = 2011
$ws.Range("E63").Value = "820"
$ws.Range("A64").Value = 430
$ws.Range("B64").Value = "Liberia"
$ws.Range("C64").Value = "GDP per Capita"
$ws.Range("D64").Value = 2012
$ws.Range("E64").Value = "853"
$ws.Range("A65").Value = 430
$ws.Range("B65").Value = "Liberia"
$ws.Range("C65").Value = "GDP per Capita"
$ws.Range("D65").Value = 2013
$ws.Range("E65").Value = "900"
$ws.Range("A66").Value = 430
$ws.Range("B66").Value = "Liberia"
$ws.Range("C66").Value = "GDP per Capita"
$ws.Range("D66").Value = 2014
$ws.Range("E66").Value = "884"
$ws.Range("A67").Value = 430
$ws.Range("B67").Value = "Liberia"
$ws.Range("C67").Value = "GDP per Capita"
$ws.Range("D67").Value = 2015
$ws.Range("E67").Value = "865"
$ws.Range("A68").Value = 430
$ws.Range("B68").Value = "Liberia"
$ws.Range("C68").Value = "GDP per Capita"
$ws.Range("D68").Value = 2016
$ws.Range("E68").Value = "829"
